# xls export geometry fix, added server power\health status
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers (row 1) to shorter labels
$ws.Range("F1").Value = "Memory tot.size"
$ws.Range("H1").Value = "Memory P/Ns"
$ws.Range("M1").Value = "HDD slot pop."
$ws.Range("N1").Value = "PSU P/Ns"

# Adjust column widths to match new, narrower headers
$ws.Columns.Item(6).ColumnWidth = 15.7109375
$ws.Columns.Item(8).ColumnWidth = 11.7109375
$ws.Columns.Item(13).ColumnWidth = 13.7109375
$ws.Columns.Item(14).ColumnWidth = 8.7109375
